$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.986.37"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.559.97"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'207.35"
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  +2.05%  "

$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("E10").Value = "  +1.84%  "

$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").Value = "1.781.96"
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").Value = "1.449.19"
$ws.Range("E13").Value = "  -6.68%  "

$ws.Range("D14").Value = "'3.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("D16").Value = "'62.08"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "26.979.48"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").Value = "0.0₃0706"
$ws.Range("E18").Value = "  +2.42%  "

$ws.Range("D19").Value = "'217.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").Value = "'7.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.21%  "

$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("D23").Value = "'9.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("D24").Value = "'1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.20%  "

$ws.Range("D25").Value = "'153.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").Value = "'15.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.18%  "

$ws.Range("E28").Value = "  +1.45%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("D30").Value = "'0.0469"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("E31").Value = "  +2.08%  "

$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("D33").Value = "1.422.72"
$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("E34").Value = "  +3.79%  "

$ws.Range("E35").Value = "  +3.47%  "

$ws.Range("D36").Value = "'1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.05%  "

$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("E38").Value = "  +0.77%  "

$ws.Range("D39").Value = "'0.531"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.71%  "

$ws.Range("E40").Value = "  +0.30%  "

$ws.Range("E42").Value = "  +0.77%  "

$ws.Range("E43").Value = "  +2.79%  "

$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.06%  "

$ws.Range("D45").Value = "'64.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.97%  "

$ws.Range("D46").Value = "'1.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("D47").Value = "1.695.71"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").Value = "'87.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.45%  "

$ws.Range("D49").Value = "'0.0521"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("D50").Value = "0.0₆01000"
$ws.Range("E50").Value = "  -0.42%  "

$ws.Range("D51").Value = "'0.0954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
